$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for the columns that change (D, J, K, L, M, O, P)
# across rows 2, 3, 4 before writing anything, since the update is a cyclic
# rotation of these values between the rows.
$cols = @("D", "J", "K", "L", "M", "O", "P")
$orig = @{}
foreach ($row in 2..4) {
    $orig[$row] = @{}
    foreach ($col in $cols) {
        $orig[$row][$col] = $ws.Range("$col$row").Value2
    }
}

# New row 2 gets old row 4's values
# New row 3 gets old row 2's values
# New row 4 gets old row 3's values
$mapping = @{ 2 = 4; 3 = 2; 4 = 3 }

foreach ($row in 2..4) {
    $srcRow = $mapping[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $orig[$srcRow][$col]
    }
}
